# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 166 (Hortaliza, Femacal de La Calera - Poroto verde),
# shifting the existing rows 166-212 down to 167-213.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(166).Insert()

$ws.Cells.Item(166, 1).Value = 3
$ws.Cells.Item(166, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(166, 3).Value = 'Coquimbo'
$ws.Cells.Item(166, 4).Value = 44463
$ws.Cells.Item(166, 5).Value = 5
$ws.Cells.Item(166, 6).Value = 100112031
$ws.Cells.Item(166, 7).Value = 'Poroto verde'
$ws.Cells.Item(166, 8).Value = 'Magnum'
$ws.Cells.Item(166, 9).Value = 'Primera'
$ws.Cells.Item(166, 10).Value = 73
$ws.Cells.Item(166, 11).Value = 33000
$ws.Cells.Item(166, 12).Value = 34000
$ws.Cells.Item(166, 13).Value = 33479
$ws.Cells.Item(166, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(166, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(166, 16).Value = 1339
$ws.Cells.Item(166, 17).Value = 25
$ws.Cells.Item(166, 18).Value = 'Hortaliza'
